$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.238.56"
$ws.Range("E2").Value = "  +0.97%  "

$ws.Range("D3").Value = "1.856.39"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4652"
$ws.Range("E7").Value = "  +0.65%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07290"
$ws.Range("E9").Value = "  -0.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8882"
$ws.Range("E10").Value = "  +1.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.09"
$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07833"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").Value = "1.839.58"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.390"
$ws.Range("E14").Value = "  +1.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.526"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.98"
$ws.Range("E16").Value = "  -0.25%  "

$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008928"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.71"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").Value = "27.266.68"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.086"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.52"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").Value = "2.077.88"
$ws.Range("E24").Value = "  -0.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.953"
$ws.Range("E25").Value = "  +5.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.75"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.39"
$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.046"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.82"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.062"
$ws.Range("E30").Value = "  -1.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08811"
$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.142"
$ws.Range("E32").Value = "  +6.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7669"
$ws.Range("E33").Value = "  +5.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.167"
$ws.Range("E34").Value = "  +3.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.510"
$ws.Range("E35").Value = "  +1.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.722"
$ws.Range("E36").Value = "  +9.99%  "

$ws.Range("E37").Value = "  +4.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01941"
$ws.Range("E38").Value = "  -0.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05213"
$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.938"
$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.028"
$ws.Range("E41").Value = "  -0.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5114"
$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1629"
$ws.Range("E43").Value = "  +0.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.447"
$ws.Range("E44").Value = "  +3.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4800"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.33"
$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9996"
$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.03"
$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.641"
$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06204"
$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.65"
$ws.Range("E51").Value = "  +1.39%  "
